$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(" methane using NIST enthalpy")

# Tout (B7): new exit/flame temperature guess
$ws.Range("B7").Value = 3000

# C7: new helper formula converting B7 to units of 1000 K
$ws.Range("C7").Formula = "=B7/1000"

# C18: moles of O2 excess term changed to 0
$ws.Range("C18").Value = 0

# C21: moles term now 1 (for CO2 branch)
$ws.Range("C21").Value = 1

# C22: moles term now 0 (N2 branch)
$ws.Range("C22").Value = 0

# J23: new integral helper formula
$ws.Range("J23").Formula = "=J21*C7+K21*C7^2/2+L21*C7^3/3+M21*C7^4/4-N21/C7+O21"

# G24: formula source changed from C22 to D21
$ws.Range("G24").Formula = "=D21*F24"

$ws.Range("B7").Select()
